# Add season-record columns (Wins, Losses, Ties) to the worksheet.
# These three new columns are appended right after the existing last
# column (AC), turning the used range from A1:AC42 into A1:AF42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers, styled like the rest of
#     the header row (bold font, thin border, centered horiz/vert). ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$header = $ws.Range("AD1:AF1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous (thin box border)

# --- Data rows (2-42): the team's season record, same for every
#     player row since it is a team-level statistic. ---
$wins = 85
$losses = 77
$ties = 0

for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins     # column AD
    $ws.Cells.Item($r, 31).Value = $losses   # column AE
    $ws.Cells.Item($r, 32).Value = $ties     # column AF
}
